$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# --- Update existing top values ---
$ws.Range("B2").Value = 16
$ws.Range("B3").Value = 560000
$ws.Range("B4").Value = 3428571.428571429

# --- Insert "Ứng lương tại CẦN THƠ" after "Công phụ phẫu 2 tại CẦN THƠ" (row 10) ---
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Ứng lương tại CẦN THƠ"
$ws.Range("B11").Value = 0

# --- Insert "Ứng lương tại LONG XUYÊN" after "Công phụ phẫu 2 tại LONG XUYÊN" (now row 18) ---
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "Ứng lương tại LONG XUYÊN"
$ws.Range("B19").Value = 0

# --- Insert "Ứng lương tại SÓC TRĂNG" after "Công phụ phẫu 2 tại SÓC TRĂNG" (now row 26) ---
$ws.Rows.Item(27).Insert()
$ws.Range("A27").Value = "Ứng lương tại SÓC TRĂNG"
$ws.Range("B27").Value = 0

# --- Append the four "Tổng lương" summary rows at the bottom (rows 28-31) ---
$ws.Range("A28").Value = "Tổng lương tại CẦN THƠ"
$ws.Range("B28").Value = 8188571.428571429

$ws.Range("A29").Value = "Tổng lương tại LONG XUYÊN"
$ws.Range("B29").Value = 0

$ws.Range("A30").Value = "Tổng lương tại SÓC TRĂNG"
$ws.Range("B30").Value = 0

$ws.Range("A31").Value = "Tổng lương"
$ws.Range("B31").Value = 8188571.428571429
